{"js": "// Update the date header and each \"two-digit \u00d7 two-digit\" answer cell.\n// Every <w:t> run in this worksheet holds a unique string, so a scoped\n// exact-text search + in-place replace (which preserves the run's\n// rFonts/sz formatting) is a safe, order-independent way to land the\n// full set of replacements described by the diff.\nconst replacements = [\n  [\"2025-12-27 Saturday\", \"2025-12-28 Sunday\"],\n  [\"63\u00d718=1134\", \"91\u00d733=3003\"],\n  [\"51\u00d773=3723\", \"94\u00d731=2914\"],\n  [\"62\u00d720=1240\", \"26\u00d725=650\"],\n  [\"38\u00d788=3344\", \"22\u00d717=374\"],\n  [\"57\u00d762=3534\", \"94\u00d749=4606\"],\n  [\"83\u00d750=4150\", \"79\u00d799=7821\"],\n  [\"30\u00d779=2370\", \"53\u00d727=1431\"],\n  [\"58\u00d720=1160\", \"36\u00d777=2772\"],\n  [\"69\u00d714=966\", \"32\u00d767=2144\"],\n  [\"45\u00d723=1035\", \"27\u00d717=459\"],\n  [\"21\u00d768=1428\", \"73\u00d795=6935\"],\n  [\"45\u00d789=4005\", \"11\u00d756=616\"],\n  [\"97\u00d757=5529\", \"52\u00d740=2080\"],\n  [\"97\u00d764=6208\", \"94\u00d750=4700\"],\n  [\"29\u00d793=2697\", \"57\u00d758=3306\"],\n  [\"72\u00d719=1368\", \"63\u00d767=4221\"],\n  [\"13\u00d747=611\", \"33\u00d785=2805\"],\n  [\"88\u00d797=8536\", \"11\u00d736=396\"],\n  [\"59\u00d779=4661\", \"12\u00d773=876\"],\n  [\"57\u00d755=3135\", \"36\u00d755=1980\"],\n  [\"95\u00d793=8835\", \"63\u00d722=1386\"],\n  [\"94\u00d774=6956\", \"88\u00d782=7216\"],\n  [\"33\u00d737=1221\", \"76\u00d760=4560\"],\n  [\"68\u00d777=5236\", \"35\u00d716=560\"],\n  [\"17\u00d718=306\", \"88\u00d792=8096\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and each \"two-digit x two-digit\" answer cell.\n# Every run in this worksheet holds a unique string, so a document-wide\n# Find/Replace (wdReplaceOne, case-sensitive) per pair is a safe,\n# order-independent way to land the full set of replacements described\n# by the diff while leaving paragraph/run formatting untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-27 Saturday\", \"2025-12-28 Sunday\"),\n    @(\"63\u00d718=1134\", \"91\u00d733=3003\"),\n    @(\"51\u00d773=3723\", \"94\u00d731=2914\"),\n    @(\"62\u00d720=1240\", \"26\u00d725=650\"),\n    @(\"38\u00d788=3344\", \"22\u00d717=374\"),\n    @(\"57\u00d762=3534\", \"94\u00d749=4606\"),\n    @(\"83\u00d750=4150\", \"79\u00d799=7821\"),\n    @(\"30\u00d779=2370\", \"53\u00d727=1431\"),\n    @(\"58\u00d720=1160\", \"36\u00d777=2772\"),\n    @(\"69\u00d714=966\", \"32\u00d767=2144\"),\n    @(\"45\u00d723=1035\", \"27\u00d717=459\"),\n    @(\"21\u00d768=1428\", \"73\u00d795=6935\"),\n    @(\"45\u00d789=4005\", \"11\u00d756=616\"),\n    @(\"97\u00d757=5529\", \"52\u00d740=2080\"),\n    @(\"97\u00d764=6208\", \"94\u00d750=4700\"),\n    @(\"29\u00d793=2697\", \"57\u00d758=3306\"),\n    @(\"72\u00d719=1368\", \"63\u00d767=4221\"),\n    @(\"13\u00d747=611\", \"33\u00d785=2805\"),\n    @(\"88\u00d797=8536\", \"11\u00d736=396\"),\n    @(\"59\u00d779=4661\", \"12\u00d773=876\"),\n    @(\"57\u00d755=3135\", \"36\u00d755=1980\"),\n    @(\"95\u00d793=8835\", \"63\u00d722=1386\"),\n    @(\"94\u00d774=6956\", \"88\u00d782=7216\"),\n    @(\"33\u00d737=1221\", \"76\u00d760=4560\"),\n    @(\"68\u00d777=5236\", \"35\u00d716=560\"),\n    @(\"17\u00d718=306\", \"88\u00d792=8096\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
